# Fruta / hortaliza, semanal
# Update Fecha/Volumen/Precio columns with refreshed weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> D, J, K, L, M, P
$data = @{
    2  = @(44839, 40, 15000, 16000, 15500, 1192)
    3  = @(44868, 30, 18000, 18000, 18000, 1385)
    4  = @(44930, 30, 17000, 17000, 17000, 1308)
    5  = @(44804, 40, 12000, 13000, 12500, 962)
    6  = @(44797, 60, 12000, 13000, 12500, 962)
    7  = @(44841, 30, 18000, 18000, 18000, 1385)
    8  = @(44943, 30, 17000, 17000, 17000, 1308)
    10 = @(44895, 30, 18000, 18000, 18000, 1385)
    11 = @(44922, 30, 17000, 17000, 17000, 1308)
    13 = @(44959, 30, 19000, 19000, 19000, 1462)
    14 = @(44874, 30, 17000, 17000, 17000, 1308)
    15 = @(44832, 60, 17000, 18000, 17500, 1346)
    16 = @(44915, 50, 18000, 18000, 18000, 1385)
    17 = @(44880, 30, 17000, 17000, 17000, 1308)
    18 = @(44894, 30, 18000, 18000, 18000, 1385)
    19 = @(44859, 30, 13000, 13000, 13000, 1000)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]  # P - Precio $/Kg
}
